$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.271060705184937
$ws.Range("B1").Value = 2.441699504852295
$ws.Range("C1").Value = 4.627000331878662
$ws.Range("D1").Value = 1.999371767044067
$ws.Range("E1").Value = 1.140095591545105
